$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.35157
$ws.Range("H2").Value = 25.05471
$ws.Range("I2").Value = 0.3629556103554933
$ws.Range("J2").Value = 0.3629556103554933
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 43.69574966666666
$ws.Range("N2").Value = 131.087249
$ws.Range("O2").Value = 0.3365063034544351
$ws.Range("P2").Value = 0.3365063034544351
$ws.Range("Q2").Value = 364.9281120436434
$ws.Range("R2").Value = 3284.353008392789
$ws.Range("S2").Value = 0.1221368507587753
$ws.Range("T2").Value = 0.1221368507587754

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.35157
$ws.Range("H3").Value = 25.05471
$ws.Range("I3").Value = 0.3629556103554933
$ws.Range("J3").Value = 0.3629556103554933
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("N3").Value = 140.44867
$ws.Range("O3").Value = 0.3605374521727266
$ws.Range("P3").Value = 0.3605374521727267
$ws.Range("Q3").Value = 390.9889663039667
$ws.Range("R3").Value = 3518.9006967357
$ws.Range("S3").Value = 0.1308590910093665
$ws.Range("T3").Value = 0.1308590910093665

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.35157
$ws.Range("H4").Value = 25.05471
$ws.Range("I4").Value = 0.3629556103554933
$ws.Range("J4").Value = 0.3629556103554933
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.08903066666667
$ws.Range("N4").Value = 54.26709200000001
$ws.Range("O4").Value = 0.1393058338430899
$ws.Range("P4").Value = 0.1393058338430899
$ws.Range("Q4").Value = 151.0718058448134
$ws.Range("R4").Value = 1359.64625260332
$ws.Range("S4").Value = 0.05056183394859962
$ws.Range("T4").Value = 0.05056183394859963

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.35157
$ws.Range("H5").Value = 25.05471
$ws.Range("I5").Value = 0.3629556103554933
$ws.Range("J5").Value = 0.3629556103554933
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.25020333333334
$ws.Range("N5").Value = 63.75061
$ws.Range("O5").Value = 0.1636504105297484
$ws.Range("P5").Value = 0.1636504105297484
$ws.Range("Q5").Value = 177.4725606525667
$ws.Range("R5").Value = 1597.2530458731
$ws.Range("S5").Value = 0.05939783463875187
$ws.Range("T5").Value = 0.05939783463875187

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.216696
$ws.Range("H6").Value = 33.650088
$ws.Range("I6").Value = 0.4874727437897329
$ws.Range("J6").Value = 0.487472743789733
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 43.69574966666666
$ws.Range("N6").Value = 131.087249
$ws.Range("O6").Value = 0.3365063034544351
$ws.Range("P6").Value = 0.3365063034544351
$ws.Range("Q6").Value = 490.1219405031013
$ws.Range("R6").Value = 4411.097464527911
$ws.Range("S6").Value = 0.1640376510474739
$ws.Range("T6").Value = 0.164037651047474

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.216696
$ws.Range("H7").Value = 33.650088
$ws.Range("I7").Value = 0.4874727437897329
$ws.Range("J7").Value = 0.487472743789733
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("N7").Value = 140.44867
$ws.Range("O7").Value = 0.3605374521727266
$ws.Range("P7").Value = 0.3605374521727267
$ws.Range("Q7").Value = 525.1233449981066
$ws.Range("R7").Value = 4726.110104982959
$ws.Range("S7").Value = 0.1757521810495987
$ws.Range("T7").Value = 0.1757521810495987

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.216696
$ws.Range("H8").Value = 33.650088
$ws.Range("I8").Value = 0.4874727437897329
$ws.Range("J8").Value = 0.487472743789733
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.08903066666667
$ws.Range("N8").Value = 54.26709200000001
$ws.Range("O8").Value = 0.1393058338430899
$ws.Range("P8").Value = 0.1393058338430899
$ws.Range("Q8").Value = 202.8991579226773
$ws.Range("R8").Value = 1826.092421304096
$ws.Range("S8").Value = 0.06790779704940765
$ws.Range("T8").Value = 0.06790779704940768

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.216696
$ws.Range("H9").Value = 33.650088
$ws.Range("I9").Value = 0.4874727437897329
$ws.Range("J9").Value = 0.487472743789733
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 21.25020333333334
$ws.Range("N9").Value = 63.75061
$ws.Range("O9").Value = 0.1636504105297484
$ws.Range("P9").Value = 0.1636504105297484
$ws.Range("Q9").Value = 238.3570707281867
$ws.Range("R9").Value = 2145.21363655368
$ws.Range("S9").Value = 0.07977511464325263
$ws.Range("T9").Value = 0.07977511464325265

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.441627666666667
$ws.Range("H10").Value = 10.324883
$ws.Range("I10").Value = 0.1495716458547737
$ws.Range("J10").Value = 0.1495716458547737
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 43.69574966666666
$ws.Range("N10").Value = 131.087249
$ws.Range("O10").Value = 0.3365063034544351
$ws.Range("P10").Value = 0.3365063034544351
$ws.Range("Q10").Value = 150.3845009685408
$ws.Range("R10").Value = 1353.460508716867
$ws.Range("S10").Value = 0.05033180164818577
$ws.Range("T10").Value = 0.05033180164818578

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.441627666666667
$ws.Range("H11").Value = 10.324883
$ws.Range("I11").Value = 0.1495716458547737
$ws.Range("J11").Value = 0.1495716458547737
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("N11").Value = 140.44867
$ws.Range("O11").Value = 0.3605374521727266
$ws.Range("P11").Value = 0.3605374521727267
$ws.Range("Q11").Value = 161.1240094728456
$ws.Range("R11").Value = 1450.11608525561
$ws.Range("S11").Value = 0.05392618011376148
$ws.Range("T11").Value = 0.05392618011376149

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.441627666666667
$ws.Range("H12").Value = 10.324883
$ws.Range("I12").Value = 0.1495716458547737
$ws.Range("J12").Value = 0.1495716458547737
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.08903066666667
$ws.Range("N12").Value = 54.26709200000001
$ws.Range("O12").Value = 0.1393058338430899
$ws.Range("P12").Value = 0.1393058338430899
$ws.Range("Q12").Value = 62.25570840558179
$ws.Range("R12").Value = 560.3013756502361
$ws.Range("S12").Value = 0.02083620284508259
$ws.Range("T12").Value = 0.02083620284508259

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.441627666666667
$ws.Range("H13").Value = 10.324883
$ws.Range("I13").Value = 0.1495716458547737
$ws.Range("J13").Value = 0.1495716458547737
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 21.25020333333334
$ws.Range("N13").Value = 63.75061
$ws.Range("O13").Value = 0.1636504105297484
$ws.Range("P13").Value = 0.1636504105297484
$ws.Range("Q13").Value = 73.13528771429223
$ws.Range("R13").Value = 658.21758942863
$ws.Range("S13").Value = 0.02447746124774385
$ws.Range("T13").Value = 0.02447746124774385
